# Insert a new data row before row 353 (shifts existing rows 353:481 down to 354:482)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(353).Insert()

# Populate the newly-inserted row 353 with its data
$ws.Range("A353").Value = 5
$ws.Range("B353").Value = 'Macroferia Regional de Talca'
$ws.Range("C353").Value = 'Maule'
$ws.Range("D353").Value = 45146
$ws.Range("E353").Value = 7
$ws.Range("F353").Value = 100112009
$ws.Range("G353").Value = 'Acelga'
$ws.Range("H353").Value = 'Sin especificar'
$ws.Range("I353").Value = 'Primera'
$ws.Range("J353").Value = 500
$ws.Range("K353").Value = 1500
$ws.Range("L353").Value = 1500
$ws.Range("M353").Value = 1500
$ws.Range("N353").Value = '$/docena de atados (4 kilos)'
$ws.Range("O353").Value = 'Región del Maule'
$ws.Range("P353").Value = 375
$ws.Range("Q353").Value = 4
$ws.Range("R353").Value = 'Hortaliza'
